$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timesheet entry: 2020-08-21 (serial 44064), 5 hours, wireframe/Adobe XD note
$ws.Range("A39").Value = 44064
$ws.Range("B39").Value = 5
$ws.Range("C39").Value = "Sivun wireframea ja adobe xd:hen tutustumista"

# Match the author's final view/selection state
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 49
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K36").Select()
